$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("L2").Value = "[4.341465385498546, 8.6634833716468]"
$ws.Range("M2").Value = [double]"7.336176555128304e-09"
$ws.Range("N2").Value = [double]"1.467235311025661e-08"
$ws.Range("P2").Value = "[-1.8365266363327724, -1.1069475616252316]"
$ws.Range("Q2").Value = [double]"2.398081733190338e-14"
$ws.Range("R2").Value = [double]"4.796163466380676e-14"
$ws.Range("T2").Value = "[7.833052100973644, 10.550421166218529]"
$ws.Range("X2").Value = [double]"4.008008008008028"
$ws.Range("Y2").Value = [double]"6.649649649649687"

# Row 3 updates
$ws.Range("L3").Value = "[4.561795289170945, 9.605192844543256]"
$ws.Range("M3").Value = [double]"7.279301872387123e-08"
$ws.Range("N3").Value = [double]"7.279301872387123e-08"
$ws.Range("P3").Value = "[1.1258159859711174, 1.9560266571900407]"
$ws.Range("Q3").Value = [double]"2.762456929872315e-12"
$ws.Range("R3").Value = [double]"2.762456929872315e-12"
$ws.Range("T3").Value = "[7.6182520369748365, 10.765868480801238]"
$ws.Range("X3").Value = [double]"16.72136136136161"
$ws.Range("Y3").Value = [double]"19.92952952952982"
